$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Bump the CodeSystem version (B3: "Version" row)
$ws.Range("B3").Value = "0.7.0"

# Remove the "Jurisdiction" / "Chile" property row entirely; everything
# below shifts up by one row and shared strings are recompacted.
$ws.Rows.Item(11).Delete()
